$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of words (German verbs) to populate column B (word),
# rows 2 through 193, replacing the previous word list in order.
$words = @(
    "wählen",
    "schwächen",
    "werfen",
    "brauen",
    "schwören",
    "zeigen",
    "sorgen",
    "spenden",
    "knarren",
    "kaufen",
    "drehen",
    "schalten",
    "strahlen",
    "bluten",
    "fischen",
    "spielen",
    "lesen",
    "zeugen",
    "kosten",
    "fließen",
    "beten",
    "tollen",
    "machen",
    "flehen",
    "pfeifen",
    "lehnen",
    "husten",
    "bauen",
    "nerven",
    "planen",
    "hassen",
    "liegen",
    "geben",
    "kichern",
    "sparen",
    "münzen",
    "freuen",
    "flüchten",
    "wecken",
    "feiern",
    "lohnen",
    "doppeln",
    "weichen",
    "ändern",
    "sichern",
    "buchen",
    "opfern",
    "tropfen",
    "sammeln",
    "ärgern",
    "wirken",
    "sprengen",
    "kürzen",
    "trennen",
    "liefern",
    "gelten",
    "spüren",
    "irren",
    "segnen",
    "deuten",
    "klingen",
    "ehren",
    "stehlen",
    "stecken",
    "danken",
    "erben",
    "sinken",
    "zünden",
    "heilen",
    "trauen",
    "werden",
    "starten",
    "bitten",
    "scheitern",
    "orten",
    "klappen",
    "formen",
    "süßen",
    "fahren",
    "wehtun",
    "suchen",
    "fällen",
    "schenken",
    "helfen",
    "kriegen",
    "jubeln",
    "rufen",
    "backen",
    "arten",
    "atmen",
    "segeln",
    "wehren",
    "graben",
    "rasen",
    "fallen",
    "meistern",
    "hören",
    "malen",
    "kehren",
    "gründen",
    "lügen",
    "lockern",
    "merken",
    "pressen",
    "haben",
    "bellen",
    "tragen",
    "enden",
    "greifen",
    "dienen",
    "zögern",
    "mögen",
    "stillen",
    "teilen",
    "scheiden",
    "warnen",
    "stammen",
    "loben",
    "baden",
    "schrecken",
    "decken",
    "gleichen",
    "schmecken",
    "fordern",
    "äußern",
    "scheinen",
    "rauchen",
    "runden",
    "töten",
    "dringen",
    "eignen",
    "hauen",
    "sterben",
    "schwingen",
    "melden",
    "boxen",
    "achten",
    "kümmern",
    "wüten",
    "jagen",
    "wiegen",
    "flüstern",
    "stoppen",
    "filmen",
    "wenden",
    "folgen",
    "foltern",
    "sperren",
    "testen",
    "rühren",
    "altern",
    "wehen",
    "zielen",
    "führen",
    "heulen",
    "fühlen",
    "grüßen",
    "schulden",
    "treiben",
    "knien",
    "schlucken",
    "wachsen",
    "schneiden",
    "tanzen",
    "füttern",
    "ziehen",
    "pflanzen",
    "reizen",
    "mauern",
    "fragen",
    "lassen",
    "fangen",
    "räumen",
    "schreiten",
    "brauchen",
    "klettern",
    "biegen",
    "träumen",
    "platzen",
    "spinnen",
    "betteln",
    "seufzen",
    "heben",
    "laufen",
    "streichen",
    "proben",
    "wundern",
    "bergen",
    "saufen",
    "siegen",
    "quälen",
    "schützen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
